$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 23).Value = "'TRUE"
}

$ws.Range("W2:W33").Select()
